# Added bugs and todos and reorganized external programs and libraries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row at row 2 (pushes the existing rows down by one) ---
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "Add search box to Edit Companies"
$ws.Range("B2").Value = "When editing a company in the Edit Companies window, it can be difficult to find the company you are looking for. It would be nice if there was a search box at the bottom."
$ws.Range("C2").Value = "Completed - Mar 1, 2010"

$ws.Range("A2:C2").Font.Bold = $false
$ws.Range("A2:C2").WrapText = $true
$ws.Range("A2:C2").VerticalAlignment = -4160
$ws.Rows.Item(2).RowHeight = 30

# --- Update the status of "Create global configs table" (now row 5) ---
$ws.Range("C5").Value = "Completed - Mar 5, 2010 (still need to go through and remove hard codes throughout program)"

# --- Append three new todo rows at the bottom of the table ---
$ws.Range("A9").Value = "Reference numbers in invoices"
$ws.Range("B9").Value = "Before the invoice system rewrite it was possible to add any number of reference numbers to an invoice such as PO, AFE, etc. This needs to be added back in."
$ws.Range("C9").Value = "OPEN"
$ws.Range("A9:C9").WrapText = $true
$ws.Range("A9:C9").VerticalAlignment = -4160
$ws.Rows.Item(9).RowHeight = 30

$ws.Range("A10").Value = "Related invoice extras to items"
$ws.Range("B10").Value = "The invoice extras section should actually be a system where extras costs can be attached to a specific line item and then printed inline with them on the final invoice. Currently they are not this way."
$ws.Range("C10").Value = "OPEN"
$ws.Range("A10:C10").WrapText = $true
$ws.Range("A10:C10").VerticalAlignment = -4160
$ws.Rows.Item(10).RowHeight = 45

$ws.Range("A11").Value = "Reorder invoice items"
$ws.Range("B11").Value = "It should be possible to reorder invoice items so that the order the are printed on the invoice is customizable. This might be done with drag and drop in the DataGrid control."
$ws.Range("C11").Value = "OPEN"
$ws.Range("A11:C11").WrapText = $true
$ws.Range("A11:C11").VerticalAlignment = -4160
$ws.Rows.Item(11).RowHeight = 30

# --- Widen columns A and C slightly to fit the new content ---
$ws.Columns.Item(1).ColumnWidth = 30.6
$ws.Columns.Item(3).ColumnWidth = 30.6

# --- Move the active selection to C2, matching the saved view state ---
$ws.Range("C2").Select() | Out-Null
